$d = $word.ActiveDocument

$pairs = @(
    @("2023-11-27 Monday", "2023-11-28 Tuesday"),
    @("47÷5=9, 2", "95÷9=10, 5"),
    @("90÷6=15, 0", "18÷8=2, 2"),
    @("45÷2=22, 1", "49÷6=8, 1"),
    @("21÷9=2, 3", "30÷8=3, 6"),
    @("65÷9=7, 2", "65÷6=10, 5"),
    @("23÷3=7, 2", "52÷6=8, 4"),
    @("25÷4=6, 1", "68÷4=17, 0"),
    @("82÷4=20, 2", "21÷5=4, 1"),
    @("61÷9=6, 7", "42÷5=8, 2"),
    @("30÷3=10, 0", "35÷2=17, 1"),
    @("69÷2=34, 1", "13÷6=2, 1"),
    @("11÷2=5, 1", "80÷9=8, 8"),
    @("46÷8=5, 6", "68÷5=13, 3"),
    @("58÷2=29, 0", "42÷7=6, 0"),
    @("99÷3=33, 0", "22÷8=2, 6"),
    @("55÷2=27, 1", "95÷9=10, 5"),
    @("53÷6=8, 5", "88÷6=14, 4"),
    @("94÷2=47, 0", "39÷6=6, 3"),
    @("36÷4=9, 0", "68÷5=13, 3"),
    @("95÷3=31, 2", "15÷9=1, 6"),
    @("59÷2=29, 1", "18÷8=2, 2"),
    @("58÷7=8, 2", "85÷2=42, 1"),
    @("69÷5=13, 4", "80÷6=13, 2"),
    @("54÷6=9, 0", "64÷7=9, 1"),
    @("59÷9=6, 5", "40÷7=5, 5")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
